$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 116.5487521575342
$ws.Range("C2").Value = 14.64110773955345
$ws.Range("D2").Value = 24.7105577116583
$ws.Range("E2").Value = 0.5115670753796644
$ws.Range("F2").Value = -32.61347595620919

$ws.Range("B3").Value = 162.8709063428561
$ws.Range("C3").Value = 18.63918410888486
$ws.Range("D3").Value = 18.08007757919806
$ws.Range("E3").Value = 0.9203049066575347
$ws.Range("F3").Value = -20.60966163179796

$ws.Range("B4").Value = 134.0590482873069
$ws.Range("C4").Value = 16.22840751237609
$ws.Range("D4").Value = 21.35107836877604
$ws.Range("E4").Value = 0.6664022896934239
$ws.Range("F4").Value = -28.97331600649987

$ws.Range("B5").Value = 113.6633471380184
$ws.Range("C5").Value = 14.36948150595614
$ws.Range("D5").Value = 20.06771767036929
$ws.Range("E5").Value = 0.616387060508636
$ws.Range("F5").Value = -33.71726392781539
